$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has: Rank (A), Team (B), ExpPoints (C).
# We need to insert 4 new columns (WIN, TOP2, TOP4, RELEGATION) between
# Team and ExpPoints, pushing ExpPoints from column C to column G.

$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).Insert()

# Fill in the new header cells (row 1)
$ws.Cells.Item(1, 3).Value = "WIN"
$ws.Cells.Item(1, 4).Value = "TOP2"
$ws.Cells.Item(1, 5).Value = "TOP4"
$ws.Cells.Item(1, 6).Value = "RELEGATION"
